$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the latest crypto snapshot.
# NumberFormat is forced to Text ("@") before assignment so values such as
# "697.32" or "7.46" are stored as literal text (matching the source feed)
# instead of being auto-converted to numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.143.33"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.859.01"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.95%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "697.32"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.62"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.855.90"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.87%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.07%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.32%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.46"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.95%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.73%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000260"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.514.89"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.855.53"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "71.206.59"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.90"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.51%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.27"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.92%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.18"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "488.24"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.722"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.50%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.61"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.63%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.50%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.53"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.24%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.74%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.012.57"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.97%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.12"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +8.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.63"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.21%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.86"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.17%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.58%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.88%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.810.05"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.93%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.26%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.52%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +12.89%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.44"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.06"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.31%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.70"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000307"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.67%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.63"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.305"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.72"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.12%  "
